$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing header cell (G1) into the new header cell (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header label
$ws.Range("H1").Value = "Save"

# Add the new data value in H2
$ws.Range("H2").Value = 0
